$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "Play Ducks Till Dawn Free | Unique Gameplay Mechanics" "Play Ducks Till Dawn for Free"
Replace-Text "Unique gameplay mechanics with moving duck symbols" "Unique gameplay mechanics with moving ducks feature"
Replace-Text "High-quality symbol and graphic design" "High-quality symbol design and easy to distinguish symbols"
Replace-Text "Engaging bonus features with free spins and multipliers" "Great attention to detail in graphics and sound design"
Replace-Text "Dark and spooky sound design adds to the game's horror theme" "Bonus features including wild and scatter symbols, free spins, and multipliers"
Replace-Text "No background music may be unappealing to some players" "No background music throughout the game"
Replace-Text "May not appeal to those who are not fans of horror themes" "Limited number of similar slots available for comparison"
Replace-Text "Read the review of Ducks Till Dawn, a spooky and engaging slot game with a unique moving duck feature. Play for free and experience the excitement." "Read our review of Ducks Till Dawn slot game and play for free. Enjoy unique gameplay mechanics and exciting bonus features."
